$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) stores values as text (inlineStr) in the workbook, not
# as numbers, so force a text number format on the cells we are about to write
# before assigning them; this stops Excel from silently re-interpreting values
# such as "245.99" or "0.00000000750" as floating point numbers and losing the
# exact original formatting/precision. (Applied range-by-range since this COM
# host does not propagate NumberFormat across every area of a union range.)
$ws.Range("D2:D22").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D40:D48").NumberFormat = "@"

# Updated coin prices (column D)
$ws.Range("D2").Value = "245.99"
$ws.Range("D3").Value = "23.98"
$ws.Range("D4").Value = "5.351"
$ws.Range("D5").Value = "0.05811"
$ws.Range("D6").Value = "3.378"
$ws.Range("D7").Value = "6.466"
$ws.Range("D8").Value = "0.8098"
$ws.Range("D9").Value = "0.9207"
$ws.Range("D10").Value = "0.1403"
$ws.Range("D11").Value = "0.07356"
$ws.Range("D12").Value = "0.03194"
$ws.Range("D13").Value = "0.03052"
$ws.Range("D14").Value = "0.09367"
$ws.Range("D15").Value = "3.844"
$ws.Range("D16").Value = "0.001558"
$ws.Range("D17").Value = "0.04718"
$ws.Range("D18").Value = "0.0005978"
$ws.Range("D19").Value = "0.006100"
$ws.Range("D20").Value = "0.001242"
$ws.Range("D21").Value = "0.004699"
$ws.Range("D22").Value = "0.00008796"
$ws.Range("D25").Value = "0.3185"
$ws.Range("D28").Value = "0.0002349"
$ws.Range("D40").Value = "0.03843"
$ws.Range("D41").Value = "0.006427"
$ws.Range("D42").Value = "0.1068"
$ws.Range("D43").Value = "0.002749"
$ws.Range("D44").Value = "0.009058"
$ws.Range("D45").Value = "0.00005246"
$ws.Range("D46").Value = "0.00000000750"
$ws.Range("D47").Value = "0.7097"
$ws.Range("D48").Value = "0.001834"

# Updated "Bestin24h" label text for rows whose 24h trend flag changed (column E)
$ws.Range("E22").Value = "21NitroExNTXBestin24h"
$ws.Range("E41").Value = "40KickTokenKICK"
